# Updated cryptos list (Price / Volume(1h) columns) with fresh quote data,
# plus a swap of the Cronos / RenderToken rows (49 and 50).
#
# Column D ("Price") is stored as text (values like "26.231.89" or
# "1.660.01" aren't valid numbers, and others need exact trailing-zero
# formatting, e.g. "3.560"), so any D cell whose new value looks like a
# plain number gets NumberFormat = "@" applied first to stop Excel from
# auto-coercing it to a numeric value and losing precision/formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.231.89"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "1.658.08"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.04"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5223"
$ws.Range("E6").Value = "  -1.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2672"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06332"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.17"
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07756"
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.439"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("D13").Value = "1.652.75"
$ws.Range("E13").Value = "  -1.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5485"
$ws.Range("E14").Value = "  -2.86%  "
$ws.Range("D15").Value = "0.0₅8263"
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.11"
$ws.Range("E16").Value = "  -1.52%  "
$ws.Range("D17").Value = "26.252.01"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.687"
$ws.Range("E19").Value = "  -3.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.58"
$ws.Range("E20").Value = "  -1.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.19"
$ws.Range("E21").Value = "  -1.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.142"
$ws.Range("E22").Value = "  -4.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.008"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "138.27"
$ws.Range("E24").Value = "  -3.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1243"
$ws.Range("E25").Value = "  -2.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.287"
$ws.Range("E26").Value = "  -2.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.13"
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("E28").Value = "  -1.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06040"
$ws.Range("E29").Value = "  -2.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.287"
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.560"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.358"
$ws.Range("E32").Value = "  -3.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.651"
$ws.Range("E33").Value = "  -3.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9837"
$ws.Range("E34").Value = "  -3.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.781"
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5960"
$ws.Range("E37").Value = "  +3.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01595"
$ws.Range("E38").Value = "  -3.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.978"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8655"
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("D41").Value = "1.043.62"
$ws.Range("E41").Value = "  -1.35%  "
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.09"
$ws.Range("D44").Value = "1.798.05"
$ws.Range("E44").Value = "  -1.83%  "
$ws.Range("D45").Value = "0.0₈109"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.38"
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.008"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.103"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05183"
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.477"
$ws.Range("E50").Value = "  +3.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4233"
$ws.Range("E51").Value = "  +0.22%  "
